$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "notes" column to Table1 (inserted after "Video Data Folder") ---
$tbl = $ws.ListObjects.Item(1)
$notesCol = $tbl.ListColumns.Add()
$notesCol.Name = "notes"
$ws.Range("W1").Value = "notes"

# --- Corrected Chamber ("M") assignments (swapped the mixed-up box groups) ---
$ws.Range("M2").Value = 11
$ws.Range("M3").Value = 12
$ws.Range("M4").Value = 5
$ws.Range("M5").Value = 6
$ws.Range("M6").Value = 13
$ws.Range("M7").Value = 14
$ws.Range("M8").Value = 3
$ws.Range("M9").Value = 4
$ws.Range("M10").Value = 11
$ws.Range("M11").Value = 12
$ws.Range("M12").Value = 5
$ws.Range("M13").Value = 6
$ws.Range("M14").Value = 13
$ws.Range("M15").Value = 14
$ws.Range("M16").Value = 3
$ws.Range("M17").Value = 4

# --- Notes / RemoveSession updates ---
$ws.Range("W9").Value = "euthanized 3/24/25 (bad implant)"
$ws.Range("S9").Value = "{0, 1}"
$ws.Range("S2").Value = "{2}"

# --- New blank, hyperlink-styled cells matching the pattern of T/U in the same rows ---
$ws.Range("V4").Style = "Hyperlink"
$ws.Range("V20").Style = "Hyperlink"

# --- Hide the raw/internal columns (Cage..LHb Target AAV) now that the sheet is cleaned up ---
$ws.Range("E1:L1").EntireColumn.Hidden = $true

# --- Re-sized visible columns (RemoveSession..notes) ---
$ws.Range("S1").EntireColumn.ColumnWidth = 15.0
$ws.Range("T1").EntireColumn.ColumnWidth = 20.666666666666668
$ws.Range("U1").EntireColumn.ColumnWidth = 22.833333333333336
$ws.Range("V1").EntireColumn.ColumnWidth = 24.666666666666668
$ws.Range("W1").EntireColumn.ColumnWidth = 62.5

# --- Selection / view bookkeeping ---
$ws.Range("S2").Select()
